$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.031763875314827
$ws.Range("D2").Value = 1.046635844979365
$ws.Range("E2").Value = 1.031284073869498
$ws.Range("F2").Value = 1.05217213072702
$ws.Range("I2").Value = 1.038758296936026
$ws.Range("J2").Value = 1.036897806767848
$ws.Range("K2").Value = 1.049400727176832
$ws.Range("L2").Value = 1.034092637259812
$ws.Range("M2").Value = 1.054921599320609
$ws.Range("N2").Value = 1.016233419032144
# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.032789164789276
$ws.Range("D3").Value = 1.047310948525577
$ws.Range("E3").Value = 1.032157712858185
$ws.Range("F3").Value = 1.05306800373172
$ws.Range("I3").Value = 1.038987252654376
$ws.Range("J3").Value = 1.03756472235329
$ws.Range("K3").Value = 1.049887868018937
$ws.Range("L3").Value = 1.034774633236877
$ws.Range("M3").Value = 1.055630051627074
$ws.Range("N3").Value = 1.016459076482328
# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.03345292540089
$ws.Range("D4").Value = 1.047748089067273
$ws.Range("E4").Value = 1.032723651785094
$ws.Range("F4").Value = 1.053648440332676
$ws.Range("I4").Value = 1.039134379162751
$ws.Range("J4").Value = 1.037996017948617
$ws.Range("K4").Value = 1.050202702644199
$ws.Range("L4").Value = 1.035215930366329
$ws.Range("M4").Value = 1.056088578221575
$ws.Range("N4").Value = 1.01660489396788
# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.033732048975581
$ws.Range("D5").Value = 1.047931934224268
$ws.Range("E5").Value = 1.032961723894424
$ws.Range("F5").Value = 1.05389263335842
$ws.Range("I5").Value = 1.039195985748246
$ws.Range("J5").Value = 1.038177275747236
$ws.Range("K5").Value = 1.050334967670114
$ws.Range("L5").Value = 1.03540145103271
$ws.Range("M5").Value = 1.056281368176288
$ws.Range("N5").Value = 1.016666147944632
# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.033778919652387
$ws.Range("D6").Value = 1.047962806774943
$ws.Range("E6").Value = 1.033001706100425
$ws.Range("F6").Value = 1.05393364481098
$ws.Range("I6").Value = 1.039206315361576
$ws.Range("J6").Value = 1.038207706259574
$ws.Range("K6").Value = 1.050357170159962
$ws.Range("L6").Value = 1.035432600711337
$ws.Range("M6").Value = 1.056313739914047
$ws.Range("N6").Value = 1.016676429951369
# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.03345665475526
$ws.Range("D7").Value = 1.047750545336916
$ws.Range("E7").Value = 1.032726832322552
$ws.Range("F7").Value = 1.053651702556123
$ws.Range("I7").Value = 1.039135203317773
$ws.Range("J7").Value = 1.037998440156261
$ws.Range("K7").Value = 1.050204470336463
$ws.Range("L7").Value = 1.03521840930594
$ws.Range("M7").Value = 1.056091154192057
$ws.Range("N7").Value = 1.016605712634294
# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.032110308212458
$ws.Range("D8").Value = 1.046863935743253
$ws.Range("E8").Value = 1.031579192012585
$ws.Range("F8").Value = 1.052474739766236
$ws.Range("I8").Value = 1.038835885169746
$ws.Range("J8").Value = 1.037123244043577
$ws.Range("K8").Value = 1.049565436423454
$ws.Range("L8").Value = 1.034323120559002
$ws.Range("M8").Value = 1.055161000537827
$ws.Range("N8").Value = 1.016309721836158
# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.029740413660409
$ws.Range("D9").Value = 1.045304000148307
$ws.Range("E9").Value = 1.02956181418017
$ws.Range("F9").Value = 1.050406559310196
$ws.Range("I9").Value = 1.038300632284336
$ws.Range("J9").Value = 1.035579199148571
$ws.Range("K9").Value = 1.048436522831822
$ws.Range("L9").Value = 1.032745538182151
$ws.Range("M9").Value = 1.053522846300444
$ws.Range("N9").Value = 1.015786645122144
# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.028162208989328
$ws.Range("D10").Value = 1.044265726253921
$ws.Range("E10").Value = 1.028220249281495
$ws.Range("F10").Value = 1.049031735034516
$ws.Range("I10").Value = 1.037938568500056
$ws.Range("J10").Value = 1.034548634837703
$ws.Range("K10").Value = 1.047682050002491
$ws.Range("L10").Value = 1.03169387623018
$ws.Range("M10").Value = 1.052431409897015
$ws.Range("N10").Value = 1.015436933976413
# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.027479240968996
$ws.Range("D11").Value = 1.04381655821374
$ws.Range("E11").Value = 1.027640142792193
$ws.Range("F11").Value = 1.048437376922608
$ws.Range("I11").Value = 1.037780556174745
$ws.Range("J11").Value = 1.034102111051802
$ws.Range("K11").Value = 1.047354924325578
$ws.Range("L11").Value = 1.031238516941321
$ws.Range("M11").Value = 1.051958977259539
$ws.Range("N11").Value = 1.015285273079877
# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.027225617222364
$ws.Range("D12").Value = 1.043649780348703
$ws.Range("E12").Value = 1.0274247863937
$ws.Range("F12").Value = 1.048216749827874
$ws.Range("I12").Value = 1.03772167802168
$ws.Range("J12").Value = 1.033936210261911
$ws.Range("K12").Value = 1.047233350978492
$ws.Range("L12").Value = 1.031069379398339
$ws.Range("M12").Value = 1.051783520798439
$ws.Range("N12").Value = 1.015228904626394
# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.027280017615007
$ws.Range("D13").Value = 1.043685551888073
$ws.Range("E13").Value = 1.027470975584802
$ws.Range("F13").Value = 1.048264068553302
$ws.Range("I13").Value = 1.037734315965915
$ws.Range("J13").Value = 1.033971798436059
$ws.Range("K13").Value = 1.047259431769757
$ws.Range("L13").Value = 1.031105659822662
$ws.Range("M13").Value = 1.051821155609231
$ws.Range("N13").Value = 1.015240997423572
# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.027458275106091
$ws.Range("D14").Value = 1.043802771007554
$ws.Range("E14").Value = 1.027622338892535
$ws.Range("F14").Value = 1.048419136873839
$ws.Range("I14").Value = 1.037775693066788
$ws.Range("J14").Value = 1.034088398494183
$ws.Range("K14").Value = 1.047344876341431
$ws.Range("L14").Value = 1.031224535911213
$ws.Range("M14").Value = 1.051944473436209
$ws.Range("N14").Value = 1.015280614355219
# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.027568113514905
$ws.Range("D15").Value = 1.043875001982907
$ws.Range("E15").Value = 1.027715614877
$ws.Range("F15").Value = 1.048514698699506
$ws.Range("I15").Value = 1.037801162316993
$ws.Range("J15").Value = 1.034160234090943
$ws.Range("K15").Value = 1.047397513067675
$ws.Range("L15").Value = 1.031297779830334
$ws.Range("M15").Value = 1.052020457103037
$ws.Range("N15").Value = 1.015305019045008
# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.028207543872253
$ws.Range("D16").Value = 1.044295544839336
$ws.Range("E16").Value = 1.028258765960842
$ws.Range("F16").Value = 1.049071200745556
$ws.Range("I16").Value = 1.037949029244233
$ws.Range("J16").Value = 1.034578263237589
$ws.Range("K16").Value = 1.047703751189398
$ws.Range("L16").Value = 1.031724097339627
$ws.Range("M16").Value = 1.052462767292265
$ws.Range("N16").Value = 1.015446994309936
# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.028608750170785
$ws.Range("D17").Value = 1.044559451351233
$ws.Range("E17").Value = 1.028599685020245
$ws.Range("F17").Value = 1.049420534976668
$ws.Range("I17").Value = 1.038041451547563
$ws.Range("J17").Value = 1.03484040655121
$ws.Range("K17").Value = 1.047895730718441
$ws.Range("L17").Value = 1.031991520081723
$ws.Range("M17").Value = 1.052740262104195
$ws.Range("N17").Value = 1.015535989269059
# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.02884280608567
$ws.Range("D18").Value = 1.044713423160159
$ws.Range("E18").Value = 1.028798614471324
$ws.Range("F18").Value = 1.04962438728749
$ws.Range("I18").Value = 1.03809524059961
$ws.Range("J18").Value = 1.034993283007494
$ws.Range("K18").Value = 1.048007667198702
$ws.Range("L18").Value = 1.032147504870935
$ws.Range("M18").Value = 1.052902136131117
$ws.Range("N18").Value = 1.015587875960449
# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.028922619726791
$ws.Range("D19").Value = 1.044765930245142
$ws.Range("E19").Value = 1.028866457332026
$ws.Range("F19").Value = 1.049693911149856
$ws.Range("I19").Value = 1.038113561026405
$ws.Range("J19").Value = 1.035045405282427
$ws.Range("K19").Value = 1.048045827506389
$ws.Range("L19").Value = 1.032200691911708
$ws.Range("M19").Value = 1.052957333673097
$ws.Range("N19").Value = 1.015605564146383
# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.028565700489789
$ws.Range("D20").Value = 1.044531132573852
$ws.Range("E20").Value = 1.028563099647718
$ws.Range("F20").Value = 1.049383045239671
$ws.Range("I20").Value = 1.038031547847356
$ws.Range("J20").Value = 1.034812283875098
$ws.Range("K20").Value = 1.047875137462503
$ws.Range("L20").Value = 1.03196282797974
$ws.Range("M20").Value = 1.052710487874797
$ws.Range("N20").Value = 1.015526443284035
# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.027405781033303
$ws.Range("D21").Value = 1.043768251144676
$ws.Range("E21").Value = 1.027577762824982
$ws.Range("F21").Value = 1.048373469133935
$ws.Range("I21").Value = 1.037763513654923
$ws.Range("J21").Value = 1.034054063838004
$ws.Range("K21").Value = 1.047319716812699
$ws.Range("L21").Value = 1.03118952977385
$ws.Range("M21").Value = 1.051908158682193
$ws.Range("N21").Value = 1.015268949116035
# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.02667684654765
$ws.Range("D22").Value = 1.04328896330404
$ws.Range("E22").Value = 1.026958942809259
$ws.Range("F22").Value = 1.047739542168882
$ws.Range("I22").Value = 1.037593917800947
$ws.Range("J22").Value = 1.033577098200291
$ws.Range("K22").Value = 1.046970130196124
$ws.Range("L22").Value = 1.030703345137119
$ws.Range("M22").Value = 1.051403853737122
$ws.Range("N22").Value = 1.015106850846655
# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.027063235011294
$ws.Range("D23").Value = 1.043543007719773
$ws.Range("E23").Value = 1.027286924356874
$ws.Range("F23").Value = 1.04807551940925
$ws.Range("I23").Value = 1.037683925261487
$ws.Range("J23").Value = 1.033829969541771
$ws.Range("K23").Value = 1.047155487602225
$ws.Range("L23").Value = 1.030961078852518
$ws.Range("M23").Value = 1.051671180646077
$ws.Range("N23").Value = 1.015192801215607
# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.028585152666763
$ws.Range("D24").Value = 1.044543928488151
$ws.Range("E24").Value = 1.028579630764475
$ws.Range("F24").Value = 1.04939998495529
$ws.Range("I24").Value = 1.038036023271862
$ws.Range("J24").Value = 1.034824991386261
$ws.Range("K24").Value = 1.047884442798075
$ws.Range("L24").Value = 1.031975792700566
$ws.Range("M24").Value = 1.052723941518101
$ws.Range("N24").Value = 1.015530756773491
# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.030352785439876
$ws.Range("D25").Value = 1.045706990102518
$ws.Range("E25").Value = 1.030082767394429
$ws.Range("F25").Value = 1.05094054076575
$ws.Range("I25").Value = 1.038439931499535
$ws.Range("J25").Value = 1.035978585983075
$ws.Range("K25").Value = 1.048728706599767
$ws.Range("L25").Value = 1.033153373385904
$ws.Range("M25").Value = 1.053946235694782
$ws.Range("N25").Value = 1.015922049337493
